$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6; this shifts existing rows 6..90 down to 7..91
# and picks up the D-column date style (s="2") from the row that follows.
$ws.Rows.Item(6).Insert()

# Fill in the new row 6 with the inserted weekly record.
$ws.Range("A6").Value = 11
$ws.Range("B6").Value = "Vega Monumental Concepción"
$ws.Range("C6").Value = "Bíobío"
$ws.Range("D6").Value = 44532
$ws.Range("E6").Value = 8
$ws.Range("F6").Value = 100112043
$ws.Range("G6").Value = "Pepino ensalada"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 250
$ws.Range("K6").Value = 6500
$ws.Range("L6").Value = 7000
$ws.Range("M6").Value = 6700
$ws.Range("N6").Value = "$/caja 60 unidades"
$ws.Range("O6").Value = "Región de Arica y Parinacota"
$ws.Range("P6").Value = 112
$ws.Range("Q6").Value = 60
$ws.Range("R6").Value = "Hortaliza"
